$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the marker_1 (column J) values for rows 5, 6, 7 to "G418"
$ws.Range("J5").Value = "G418"
$ws.Range("J6").Value = "G418"
$ws.Range("J7").Value = "G418"

# Update the active cell selection to J7 (matches the diff's sheetView selection)
$ws.Range("J7").Select()
